$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A144").Value = "IMX-USD"
$ws.Range("A145").Value = "GRT-USD"
